$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet1 "Recommandations": refreshed BRVM stats, re-sorted by Variation Totale (%) desc ---
# row 2
$ws1.Range("A2").Value = "NEI-CEDA CI"
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 3
$ws1.Range("D2").Value = 2805
$ws1.Range("E2").Value = 950
$ws1.Range("F2").Value = "🟡 Observer"
$ws1.Range("G2").Value = "➖ Neutre"
# row 3
$ws1.Range("A3").Value = "BRVM - SERVICES PUBLICS"
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 6
$ws1.Range("D3").Value = 2568
$ws1.Range("E3").Value = 111.96
$ws1.Range("F3").Value = "🟡 Observer"
$ws1.Range("G3").Value = "➖ Neutre"
# row 4
$ws1.Range("A4").Value = "SUCRIVOIRE"
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 2
$ws1.Range("D4").Value = 1985
$ws1.Range("E4").Value = 995
$ws1.Range("F4").Value = "🟡 Observer"
$ws1.Range("G4").Value = "➖ Neutre"
# row 5
$ws1.Range("A5").Value = "BRVM - AUTRES SECTEURS"
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 3
$ws1.Range("D5").Value = 1792.1
$ws1.Range("E5").Value = 585.48
$ws1.Range("F5").Value = "🟡 Observer"
$ws1.Range("G5").Value = "➖ Neutre"
# row 6
$ws1.Range("A6").Value = "BRVM - DISTRIBUTION"
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 3
$ws1.Range("D6").Value = 1616.32
$ws1.Range("E6").Value = 536.49
$ws1.Range("F6").Value = "🟡 Observer"
$ws1.Range("G6").Value = "➖ Neutre"
# row 7
$ws1.Range("A7").Value = "AIR LIQUIDE CI"
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 2
$ws1.Range("D7").Value = 1405
$ws1.Range("E7").Value = 705
$ws1.Range("F7").Value = "🟡 Observer"
$ws1.Range("G7").Value = "➖ Neutre"
# row 8
$ws1.Range("A8").Value = "BRVM - TRANSPORT"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 3
$ws1.Range("D8").Value = 1053.75
$ws1.Range("E8").Value = 342.68
$ws1.Range("F8").Value = "🟡 Observer"
$ws1.Range("G8").Value = "➖ Neutre"
# row 9
$ws1.Range("A9").Value = "BRVM - AGRICULTURE"
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 3
$ws1.Range("D9").Value = 1015.52
$ws1.Range("E9").Value = 338.92
$ws1.Range("F9").Value = "🟡 Observer"
$ws1.Range("G9").Value = "➖ Neutre"
# row 10
$ws1.Range("A10").Value = "ERIUM CI"
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 1
$ws1.Range("D10").Value = 810
$ws1.Range("E10").Value = 810
$ws1.Range("F10").Value = "🟡 Observer"
$ws1.Range("G10").Value = "➖ Neutre"
# row 11
$ws1.Range("A11").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 3
$ws1.Range("D11").Value = 558.69
$ws1.Range("E11").Value = 183.58
$ws1.Range("F11").Value = "🟡 Observer"
$ws1.Range("G11").Value = "➖ Neutre"
# row 12
$ws1.Range("A12").Value = "BRVM - INDUSTRIE                (**)"
$ws1.Range("B12").Value = 0
$ws1.Range("C12").Value = 2
$ws1.Range("D12").Value = 538.18
$ws1.Range("E12").Value = 269.26
$ws1.Range("F12").Value = "🟡 Observer"
$ws1.Range("G12").Value = "➖ Neutre"
# row 13
$ws1.Range("A13").Value = "BRVM - FINANCES"
$ws1.Range("B13").Value = 0
$ws1.Range("C13").Value = 3
$ws1.Range("D13").Value = 449.67
$ws1.Range("E13").Value = 150.01
$ws1.Range("F13").Value = "🟡 Observer"
$ws1.Range("G13").Value = "➖ Neutre"
# row 14
$ws1.Range("A14").Value = "BRVM-PRINCIPAL                   (**)"
$ws1.Range("B14").Value = 0
$ws1.Range("C14").Value = 2
$ws1.Range("D14").Value = 444.8
$ws1.Range("E14").Value = 222.42
$ws1.Range("F14").Value = "🟡 Observer"
$ws1.Range("G14").Value = "➖ Neutre"
# row 15
$ws1.Range("A15").Value = "BRVM - SERVICES FINANCIERS"
$ws1.Range("B15").Value = 0
$ws1.Range("C15").Value = 3
$ws1.Range("D15").Value = 441.93
$ws1.Range("E15").Value = 147.43
$ws1.Range("F15").Value = "🟡 Observer"
$ws1.Range("G15").Value = "➖ Neutre"
# row 16
$ws1.Range("A16").Value = "BRVM-PRESTIGE"
$ws1.Range("B16").Value = 0
$ws1.Range("C16").Value = 3
$ws1.Range("D16").Value = 440.87
$ws1.Range("E16").Value = 146.76
$ws1.Range("F16").Value = "🟡 Observer"
$ws1.Range("G16").Value = "➖ Neutre"
# row 17
$ws1.Range("A17").Value = "BRVM - INDUSTRIELS"
$ws1.Range("B17").Value = 0
$ws1.Range("C17").Value = 3
$ws1.Range("D17").Value = 368.46
$ws1.Range("E17").Value = 121.78
$ws1.Range("F17").Value = "🟡 Observer"
$ws1.Range("G17").Value = "➖ Neutre"
# row 18
$ws1.Range("A18").Value = "BRVM - ENERGIE"
$ws1.Range("B18").Value = 0
$ws1.Range("C18").Value = 3
$ws1.Range("D18").Value = 342.14
$ws1.Range("E18").Value = 114.81
$ws1.Range("F18").Value = "🟡 Observer"
$ws1.Range("G18").Value = "➖ Neutre"
# row 19
$ws1.Range("A19").Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Range("B19").Value = 0
$ws1.Range("C19").Value = 3
$ws1.Range("D19").Value = 290.7
$ws1.Range("E19").Value = 96.77
$ws1.Range("F19").Value = "🟡 Observer"
$ws1.Range("G19").Value = "➖ Neutre"
# row 20
$ws1.Range("A20").Value = "BRVM - CONSOMMATION DE BASE         (**)"
$ws1.Range("B20").Value = 0
$ws1.Range("C20").Value = 1
$ws1.Range("D20").Value = 224.89
$ws1.Range("E20").Value = 224.89
$ws1.Range("F20").Value = "🟡 Observer"
$ws1.Range("G20").Value = "➖ Neutre"
# row 21
$ws1.Range("A21").Value = "BRVM - CONSOMMATION DE BASE          (**)"
$ws1.Range("B21").Value = 0
$ws1.Range("C21").Value = 1
$ws1.Range("D21").Value = 224.67
$ws1.Range("E21").Value = 224.67
$ws1.Range("F21").Value = "🟡 Observer"
$ws1.Range("G21").Value = "➖ Neutre"
# row 22
$ws1.Range("A22").Value = "SAFCA CI (SAFC)"
$ws1.Range("B22").Value = 2
$ws1.Range("C22").Value = 0
$ws1.Range("D22").Value = 10.87
$ws1.Range("E22").Value = 7.38
$ws1.Range("F22").Value = "🟡 Observer"
$ws1.Range("G22").Value = "➖ Neutre"
# row 23
$ws1.Range("A23").Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Range("B23").Value = 2
$ws1.Range("C23").Value = 1
$ws1.Range("D23").Value = 9.04
$ws1.Range("E23").Value = 3.27
$ws1.Range("F23").Value = "🟡 Observer"
$ws1.Range("G23").Value = "👀 À surveiller"
# row 24
$ws1.Range("A24").Value = "ERIUM CI (SIVC)"
$ws1.Range("B24").Value = 1
$ws1.Range("C24").Value = 0
$ws1.Range("D24").Value = 7.28
$ws1.Range("E24").Value = 7.28
$ws1.Range("F24").Value = "🟡 Observer"
$ws1.Range("G24").Value = "➖ Neutre"
# row 25
$ws1.Range("A25").Value = "AIR LIQUIDE CI (SIVC)"
$ws1.Range("B25").Value = 1
$ws1.Range("C25").Value = 0
$ws1.Range("D25").Value = 7.09
$ws1.Range("E25").Value = 7.09
$ws1.Range("F25").Value = "🟡 Observer"
$ws1.Range("G25").Value = "➖ Neutre"
# row 26
$ws1.Range("A26").Value = "SETAO CI (STAC)"
$ws1.Range("B26").Value = 1
$ws1.Range("C26").Value = 0
$ws1.Range("D26").Value = 6.78
$ws1.Range("E26").Value = 6.78
$ws1.Range("F26").Value = "🟡 Observer"
$ws1.Range("G26").Value = "➖ Neutre"
# row 27
$ws1.Range("A27").Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 0
$ws1.Range("D27").Value = 5.26
$ws1.Range("E27").Value = 5.26
$ws1.Range("F27").Value = "🟡 Observer"
$ws1.Range("G27").Value = "➖ Neutre"
# row 28
$ws1.Range("A28").Value = "BERNABE CI (BNBC)"
$ws1.Range("B28").Value = 1
$ws1.Range("C28").Value = 1
$ws1.Range("D28").Value = 3.96
$ws1.Range("E28").Value = -3.45
$ws1.Range("F28").Value = "🟡 Observer"
$ws1.Range("G28").Value = "👀 À surveiller"
# row 29
$ws1.Range("A29").Value = "BICI CI (BICC)"
$ws1.Range("B29").Value = 1
$ws1.Range("C29").Value = 0
$ws1.Range("D29").Value = 2.54
$ws1.Range("E29").Value = 2.54
$ws1.Range("F29").Value = "🟡 Observer"
$ws1.Range("G29").Value = "➖ Neutre"
# row 30
$ws1.Range("A30").Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Range("B30").Value = 1
$ws1.Range("C30").Value = 1
$ws1.Range("D30").Value = 0.83
$ws1.Range("E30").Value = -2.5
$ws1.Range("F30").Value = "🟡 Observer"
$ws1.Range("G30").Value = "👀 À surveiller"
# row 31
$ws1.Range("A31").Value = "ORAGROUP TOGO (ORGT)"
$ws1.Range("B31").Value = 1
$ws1.Range("C31").Value = 1
$ws1.Range("D31").Value = 0.13
$ws1.Range("E31").Value = -3.47
$ws1.Range("F31").Value = "🟡 Observer"
$ws1.Range("G31").Value = "👀 À surveiller"
# row 32
$ws1.Range("A32").Value = "TOTAL"
$ws1.Range("B32").Value = 0
$ws1.Range("C32").Value = 3
$ws1.Range("D32").Value = 0
$ws1.Range("E32").Value = 0
$ws1.Range("F32").Value = "🟡 Observer"
$ws1.Range("G32").Value = "➖ Neutre"
# row 33
$ws1.Range("A33").Value = "SMB CI (SMBC)"
$ws1.Range("B33").Value = 1
$ws1.Range("C33").Value = 1
$ws1.Range("D33").Value = -0.07
$ws1.Range("E33").Value = 2.89
$ws1.Range("F33").Value = "🟡 Observer"
$ws1.Range("G33").Value = "👀 À surveiller"
# row 34
$ws1.Range("A34").Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Range("B34").Value = 1
$ws1.Range("C34").Value = 1
$ws1.Range("D34").Value = -0.48
$ws1.Range("E34").Value = 3.13
$ws1.Range("F34").Value = "🟡 Observer"
$ws1.Range("G34").Value = "👀 À surveiller"
# row 35
$ws1.Range("A35").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("B35").Value = 0
$ws1.Range("C35").Value = 1
$ws1.Range("D35").Value = -1.14
$ws1.Range("E35").Value = -1.14
$ws1.Range("F35").Value = "🟡 Observer"
$ws1.Range("G35").Value = "➖ Neutre"
# row 36
$ws1.Range("A36").Value = "PALM CI (PALC)"
$ws1.Range("B36").Value = 0
$ws1.Range("C36").Value = 1
$ws1.Range("D36").Value = -1.54
$ws1.Range("E36").Value = -1.54
$ws1.Range("F36").Value = "🟡 Observer"
$ws1.Range("G36").Value = "➖ Neutre"
# row 37
$ws1.Range("A37").Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws1.Range("B37").Value = 0
$ws1.Range("C37").Value = 1
$ws1.Range("D37").Value = -1.91
$ws1.Range("E37").Value = -1.91
$ws1.Range("F37").Value = "🟡 Observer"
$ws1.Range("G37").Value = "➖ Neutre"
# row 38
$ws1.Range("A38").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("B38").Value = 0
$ws1.Range("C38").Value = 1
$ws1.Range("D38").Value = -2.06
$ws1.Range("E38").Value = -2.06
$ws1.Range("F38").Value = "🟡 Observer"
$ws1.Range("G38").Value = "➖ Neutre"
# row 39
$ws1.Range("A39").Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$ws1.Range("B39").Value = 0
$ws1.Range("C39").Value = 1
$ws1.Range("D39").Value = -2.78
$ws1.Range("E39").Value = -2.78
$ws1.Range("F39").Value = "🟡 Observer"
$ws1.Range("G39").Value = "➖ Neutre"
# row 40
$ws1.Range("A40").Value = "UNIWAX CI (UNXC)"
$ws1.Range("B40").Value = 0
$ws1.Range("C40").Value = 1
$ws1.Range("D40").Value = -3.11
$ws1.Range("E40").Value = -3.11
$ws1.Range("F40").Value = "🟡 Observer"
$ws1.Range("G40").Value = "➖ Neutre"
# row 41
$ws1.Range("A41").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Range("B41").Value = 0
$ws1.Range("C41").Value = 1
$ws1.Range("D41").Value = -3.45
$ws1.Range("E41").Value = -3.45
$ws1.Range("F41").Value = "🟡 Observer"
$ws1.Range("G41").Value = "➖ Neutre"
# row 42
$ws1.Range("A42").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Range("B42").Value = 0
$ws1.Range("C42").Value = 1
$ws1.Range("D42").Value = -4.76
$ws1.Range("E42").Value = -4.76
$ws1.Range("F42").Value = "🟡 Observer"
$ws1.Range("G42").Value = "➖ Neutre"

# --- Clear now-unused rows 43-46 left over from the previous (longer) sorted list ---
$ws1.Range("A43:G46").ClearContents()

# --- Sheet2 "Top_YTD": refreshed YTD progression figures, re-sorted ---
# row 2
$ws2.Range("A2").Value = "BRVM - SERVICES PUBLICS"
$ws2.Range("B2").Value = 572498.37
# row 3
$ws2.Range("A3").Value = "NEI-CEDA CI"
$ws2.Range("B3").Value = 110701.25
# row 4
$ws2.Range("A4").Value = "BRVM - AUTRES SECTEURS"
$ws2.Range("B4").Value = 33803.23
# row 5
$ws2.Range("A5").Value = "BRVM - DISTRIBUTION"
$ws2.Range("B5").Value = 25963.7
# row 6
$ws2.Range("A6").Value = "SUCRIVOIRE"
$ws2.Range("B6").Value = 11835.5
# row 7
$ws2.Range("A7").Value = "BRVM - TRANSPORT"
$ws2.Range("B7").Value = 9086.13
# row 8
$ws2.Range("A8").Value = "BRVM - AGRICULTURE"
$ws2.Range("B8").Value = 8331.92
# row 9
$ws2.Range("A9").Value = "AIR LIQUIDE CI"
$ws2.Range("B9").Value = 6340
# row 10
$ws2.Range("A10").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Range("B10").Value = 2244.86
# row 11
$ws2.Range("A11").Value = "BRVM - FINANCES"
$ws2.Range("B11").Value = 1460.44

Write-Host "Edit complete."